$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B27").Value = 6556
$ws.Range("C27").Value = 1017
$ws.Range("D27").Value = 6115050
$ws.Range("E27").Value = 932.7410006101281
$ws.Range("F27").Value = 10.18487394957983
$ws.Range("G27").Value = 7.391763463569156
$ws.Range("H27").Value = 25.42387899350775
